$d = $word.ActiveDocument

# --- First paragraph: replace the ID placeholder text and drop the
#     trailing whitespace run, keeping the paragraph mark intact.
$p1 = $d.Paragraphs(1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End
$body = $d.Range($p1Start, $p1End - 1)
$body.Text = "**ID__AFFARS_COVID__ID**"

# --- First paragraph formatting: indent + paragraph border (5-twip
#     spacing on all four edges, no explicit border line).
$p1.Format.LeftIndent = 11.25
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromRight = 5
